$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Add the new shared string used by the inserted row ("Force").
# (No direct "add shared string" API - it's created implicitly the first
# time a cell's .Value is set to this text.)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Insert a new row above row 8 by manually shifting rows 8..17 down to
# 9..18 (bottom-up so nothing is overwritten before it's copied). We copy
# cell content+format for columns A:G only, so we never touch the full
# 16384-column row (which some bulk "Insert"/whole-row Copy operations do
# in this engine, bloating the sheet).
# ---------------------------------------------------------------------
for ($r = 17; $r -ge 8; $r--) {
    $ws.Range("A" + $r + ":G" + $r).Copy($ws.Range("A" + ($r + 1) + ":G" + ($r + 1)))
}

# Range.Copy only moves cell-level formatting; row-level attributes (row
# height, row default style) stay put on their original row number, so
# fix those up to follow the content that now lives there.

# Row 9 now holds what used to be row 8 (the boolean-settings row) -
# give it row 8's old row-level look: default height, style 3 on the row.
$ws.Rows.Item(9).AutoFit()

# Row 10 now holds what used to be row 9 (the column-header row) - it
# needs the taller 27pt row height that travels with that content.
$ws.Rows.Item(10).RowHeight = 27

# ---------------------------------------------------------------------
# Populate the newly freed-up row 8 with the "Force" toggle row, cloning
# the look of the boolean rows above it (style index 3: bold 11pt 宋体,
# orange fill, thin box border, left/center aligned, wrapped).
# ---------------------------------------------------------------------
$row8 = $ws.Rows.Item(8)
$row8.Font.Name = "宋体"
$row8.Font.Size = 11
$row8.Font.Bold = $true
$row8.Interior.Pattern = 1      # xlSolid
$row8.Interior.Color = 49407    # RGB(255,192,0) -> matches fill FFFFC000
$row8.VerticalAlignment = -4108 # xlCenter
$row8.WrapText = $true
$row8.Borders.LineStyle = 1     # xlContinuous
$row8.Borders.Weight = 2        # xlThin

$rngA8G8 = $ws.Range("A8:G8")
$rngA8G8.HorizontalAlignment = -4131  # xlLeft

$ws.Cells.Item(8, 1).Value = "Force"
$ws.Cells.Item(8, 2).Value = $false
$ws.Cells.Item(8, 3).Value = $false
$ws.Cells.Item(8, 4).Value = $false
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = $false
$ws.Cells.Item(8, 7).Value = $false

# ---------------------------------------------------------------------
# View / selection bookkeeping so the saved file matches what Excel would
# have written after this edit: the frozen pane now splits after the
# (one row taller) header block and the active cell moved to A9.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A11").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A9").Select()

# ---------------------------------------------------------------------
# Data validation ranges grow by one row to include the new row 8.
# ---------------------------------------------------------------------
$ws.Range("A6:A8").Validation.Delete()
$ws.Range("A6:A9").Validation.Add(1, 1, 1)

$ws.Range("B6:G8").Validation.Delete()
$ws.Range("B6:G9").Validation.Add(3, 1, 2, "TRUE,FALSE")
